# Update xlbean version to v0.3.0
# - D1 header changes from "list#value?toMap=value&type=string" to
#   "list#value?toMap=value&readAs=text"
# - Row 3 (previously blank) is removed, shifting rows 4-8 up to rows 3-7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("toBean")

# Update the header text in D1 to use the new 'readAs' query param
$ws.Range("D1").Value2 = "list#value?toMap=value&readAs=text"

# Delete the empty row 3 so that the data below shifts up by one row
$ws.Rows("3").Delete()
